$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 10
$ws.Range("B3").Value = 5.5
$ws.Range("C4").Value = 1.4

# Update column widths: columns A and C get new, wider explicit widths
# (replacing the old bestFit auto-sized widths); column B goes back to
# the workbook's standard/default column width (no longer a bestFit
# override sized for its old, now-removed, custom width).
#
# NOTE: Excel's ColumnWidth property is expressed in characters of the
# Normal style font, while the width actually persisted to the sheet's
# <col> XML is that value re-measured on a pixel grid (padded +5px then
# divided back down), so the saved width is always a little larger than
# what was assigned. Subtract that fixed padding (5 / max-digit-width)
# before assigning so the exported width lands on the desired value.
$mdw = 7
$pad = 5 / $mdw
$ws.Columns.Item(1).ColumnWidth = 27 - $pad
$ws.Columns.Item(3).ColumnWidth = 27.25 - $pad
$ws.Columns.Item(2).ColumnWidth = $ws.StandardWidth
